$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Arinj Coliving"
$ws.Range("C16").Value = "Coliving"
$ws.Range("D16").Value = 30
$ws.Range("E16").Value = 40.2323582770568
$ws.Range("F16").Value = 44.5704503485139
$ws.Range("G16").Value = "Avan"

$ws.Range("B20").Select()
